# feat: working push products demo
#
# Remove the "brandName" column (column X) from the product feed sheet —
# every cell to its right (productTags_enUS/svSE/itIT, formerly Y/Z/AA)
# shifts one column to the left, taking its values/styles/widths with it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(24).Delete()

# Product group "2" (Sweatpants) now also references the parent group SKU
# alongside its own SKU. Re-apply the cell's text-quote-prefix formatting
# (shared with the rest of column C) after the write, since assigning
# .Value resets the cell to the default style.
$ws.Range("C6").Value = "1-001,1-001-002"
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Leave the selection where the author left it when they saved.
$ws.Range("C7").Select()
